$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change existing Gender value from "f" to "m" (shared string is shared across all cells using it)
$ws.Range("B2").Value = "m"
$ws.Range("B3").Value = "m"

# Add new rows 4 and 5
$ws.Range("A4").Value = "Adidas"
$ws.Range("B4").Value = "m"
$ws.Range("C4").Value = 37
$ws.Range("D4").Value = 7
$ws.Range("E4").Value = 5.5
$ws.Range("F4").Formula = '="38 1/3"'

$ws.Range("A5").Value = "Adidas"
$ws.Range("B5").Value = "m"
$ws.Range("C5").Value = 37.5
$ws.Range("D5").Value = 7.5
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 39

$ws.Range("F6").Select()
